{"js": "// Update the worksheet date and every division problem's operands, per the\n// commit's regenerated numbers. Every \"old\" string below occurs exactly once\n// in the document, so a plain matchCase search+replace is unambiguous.\nconst replacements = [\n  [\"2024-07-21 Sunday\", \"2024-07-22 Monday\"],\n  [\"92\u00f72=\", \"41\u00f72=\"],\n  [\"25\u00f74=\", \"14\u00f77=\"],\n  [\"36\u00f73=\", \"27\u00f78=\"],\n  [\"93\u00f75=\", \"33\u00f74=\"],\n  [\"47\u00f76=\", \"48\u00f74=\"],\n  [\"52\u00f73=\", \"51\u00f79=\"],\n  [\"65\u00f76=\", \"88\u00f72=\"],\n  [\"97\u00f74=\", \"30\u00f72=\"],\n  [\"40\u00f76=\", \"36\u00f77=\"],\n  [\"45\u00f72=\", \"27\u00f75=\"],\n  [\"53\u00f75=\", \"33\u00f79=\"],\n  [\"13\u00f79=\", \"30\u00f79=\"],\n  [\"14\u00f76=\", \"58\u00f73=\"],\n  [\"32\u00f75=\", \"75\u00f72=\"],\n  [\"52\u00f76=\", \"51\u00f73=\"],\n  [\"50\u00f72=\", \"25\u00f77=\"],\n  [\"78\u00f73=\", \"12\u00f75=\"],\n  [\"75\u00f74=\", \"62\u00f74=\"],\n  [\"45\u00f75=\", \"22\u00f75=\"],\n  [\"85\u00f73=\", \"39\u00f72=\"],\n  [\"18\u00f78=\", \"30\u00f78=\"],\n  [\"22\u00f79=\", \"31\u00f79=\"],\n  [\"24\u00f76=\", \"77\u00f74=\"],\n  [\"10\u00f76=\", \"20\u00f77=\"],\n  [\"62\u00f78=\", \"81\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every division problem's operands, per the\n# commit's regenerated numbers. Every \"old\" string below occurs exactly once\n# in the document, so Find/Replace on the whole body is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-21 Sunday\", \"2024-07-22 Monday\"),\n    @(\"92\u00f72=\", \"41\u00f72=\"),\n    @(\"25\u00f74=\", \"14\u00f77=\"),\n    @(\"36\u00f73=\", \"27\u00f78=\"),\n    @(\"93\u00f75=\", \"33\u00f74=\"),\n    @(\"47\u00f76=\", \"48\u00f74=\"),\n    @(\"52\u00f73=\", \"51\u00f79=\"),\n    @(\"65\u00f76=\", \"88\u00f72=\"),\n    @(\"97\u00f74=\", \"30\u00f72=\"),\n    @(\"40\u00f76=\", \"36\u00f77=\"),\n    @(\"45\u00f72=\", \"27\u00f75=\"),\n    @(\"53\u00f75=\", \"33\u00f79=\"),\n    @(\"13\u00f79=\", \"30\u00f79=\"),\n    @(\"14\u00f76=\", \"58\u00f73=\"),\n    @(\"32\u00f75=\", \"75\u00f72=\"),\n    @(\"52\u00f76=\", \"51\u00f73=\"),\n    @(\"50\u00f72=\", \"25\u00f77=\"),\n    @(\"78\u00f73=\", \"12\u00f75=\"),\n    @(\"75\u00f74=\", \"62\u00f74=\"),\n    @(\"45\u00f75=\", \"22\u00f75=\"),\n    @(\"85\u00f73=\", \"39\u00f72=\"),\n    @(\"18\u00f78=\", \"30\u00f78=\"),\n    @(\"22\u00f79=\", \"31\u00f79=\"),\n    @(\"24\u00f76=\", \"77\u00f74=\"),\n    @(\"10\u00f76=\", \"20\u00f77=\"),\n    @(\"62\u00f78=\", \"81\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
